# "new comments on database"
#
# Target edit (per the commit's xlsx diff):
#  - Grupos (sheet1): add a new column H, value 1, to every data row (1..156).
#    The old "recomendations" threaded comment that lived on H1 is removed,
#    and the "availibility" threaded comment that lived on I1 moves onto the
#    new H1 (which now also carries the literal value 1).
#  - Cursos (sheet2): two new threaded comments, "recomendations" on D1 and
#    "description" on E1 (header-documentation only, no data written).
#  - Profesores (sheet3): one new threaded comment, "bio" on C1.
#  - Áreas (sheet4): rows 5 and 6 in column A swap values
#    ("Capacitación" <-> "Arte y Salud").
#  - View-state: Grupos becomes the active/selected sheet (was Profesores),
#    with per-sheet selections updated to match the saved file.

$wb = $excel.ActiveWorkbook

$wsGrupos     = $wb.Worksheets.Item(1)   # Grupos
$wsCursos     = $wb.Worksheets.Item(2)   # Cursos
$wsProfesores = $wb.Worksheets.Item(3)   # Profesores
$wsAreas      = $wb.Worksheets.Item(4)   # Áreas

# --- Grupos: drop the old H1 "recomendations" comment, move the I1
#     "availibility" comment onto H1, then populate the new H column. ---
$wsGrupos.Range("H1").Comment.Delete()

$availabilityText = $wsGrupos.Range("I1").Comment.Text()
$wsGrupos.Range("I1").Comment.Delete()
$wsGrupos.Range("H1").AddCommentThreaded($availabilityText)

$wsGrupos.Range("H1:H156").Value = 1

# --- Cursos: two brand-new threaded comments documenting extra columns. ---
$wsCursos.Range("D1").AddCommentThreaded("recomendations")
$wsCursos.Range("E1").AddCommentThreaded("description")

# --- Profesores: one brand-new threaded comment. ---
$wsProfesores.Range("C1").AddCommentThreaded("bio")

# --- Áreas: swap the two area names in A5/A6. ---
$wsAreas.Range("A5").Value = "Arte y Salud"
$wsAreas.Range("A6").Value = "Capacitación"

# --- View state: restore per-sheet selections, finishing on Grupos so it
#     ends up as the active/saved tab (matching the new tabSelected owner). ---
$wsCursos.Activate()
$wsCursos.Range("F10").Select()

$wsProfesores.Activate()
$wsProfesores.Range("F3").Select()

$wsAreas.Activate()
$wsAreas.Range("B6").Select()

$wsGrupos.Activate()
$wsGrupos.Range("I153").Select()
